$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the H1 title.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Replace the final "Prompt: ..." image-generation paragraph's italic text
#    with the meta description sentence (keeps the run's italic formatting).
$d.Content.Find.Execute(
    "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Wild Pumpkins" + [char]34 + `
    ". The image should feature a happy Maya warrior with glasses. The warrior should be holding a Jack O'Lantern " + `
    "and standing in the middle of a spooky cemetery with tombstones, bats, and fog in the background. The image " + `
    "should convey a Halloween and gothic atmosphere. Use bold colors and intricate details to grab the attention " + `
    "of potential players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Wild Pumpkins slot game, with its Halloween theme, bonus features, and 4096 ways to win. Play for free!",
    2)

# 3. Insert a new bold paragraph ("Play Wild Pumpkins Free: Review, RTP, Features")
#    right before that paragraph, matching the structure of the other runs
#    (leading empty run + bold run). The insertion point is placed right
#    before the previous paragraph's mark (End-1) rather than at the start
#    of the following paragraph (End), because InsertXML at a range
#    collapsed to the very start of a paragraph replaces that paragraph's
#    content instead of inserting a new paragraph ahead of it.
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($count - 1)
$insertPos = $prevPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Wild Pumpkins Free: Review, RTP, Features</w:t></w:r></w:p>"
$insertRange.InsertXML($newParaXml)
